# Update countries & provincias Spain
# Data refresh: 23 de Marzo de 2020, 09:46 -> 10:16
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 10:16"

$ws.Range("A12").Value = "Suiza"
$ws.Range("B12").Value = 7776
$ws.Range("C12").Value = 302
$ws.Range("D12").Value = 131
$ws.Range("E12").Value = 7545
$ws.Range("F12").Value = 141
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 100

$ws.Range("A18").Value = "Suecia"
$ws.Range("B18").Value = 1934
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 16
$ws.Range("E18").Value = 1894
$ws.Range("F18").Value = 76
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 24

$ws.Range("A30").Value = "Pakistan"
$ws.Range("B30").Value = 804
$ws.Range("C30").Value = 28
$ws.Range("D30").Value = 13
$ws.Range("E30").Value = 785
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 6

$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 579
$ws.Range("C39").Value = 65
$ws.Range("D39").Value = 30
$ws.Range("E39").Value = 500
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 49

$ws.Range("A40").Value = "Islandia"
$ws.Range("B40").Value = 568
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 5
$ws.Range("E40").Value = 562
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 1

$ws.Range("A48").Value = "Eslovenia"
$ws.Range("B48").Value = 414
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 411
$ws.Range("F48").Value = 12
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 3

$ws.Range("A88").Value = "Albania"
$ws.Range("B88").Value = 89
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 2
$ws.Range("E88").Value = 83
$ws.Range("F88").Value = 2
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 4

$ws.Range("A100").Value = "Kazajistan"
$ws.Range("B100").Value = 62
$ws.Range("C100").Value = 2
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 62
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 0

$ws.Range("A112").Value = "Nigeria"
$ws.Range("B112").Value = 35
$ws.Range("C112").Value = 5
$ws.Range("D112").Value = 2
$ws.Range("E112").Value = 33
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 0

$ws.Range("A113").Value = "Consejo Danes para los Refugiados"
$ws.Range("B113").Value = 30
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 0
$ws.Range("E113").Value = 29
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 1

$ws.Range("A120").Value = "Ghana"
$ws.Range("B120").Value = 23
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 0
$ws.Range("E120").Value = 22
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 1

$ws.Range("A121").Value = "Monaco"
$ws.Range("B121").Value = 23
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 1
$ws.Range("E121").Value = 22
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 0

$ws.Range("A122").Value = "Puerto Rico"
$ws.Range("B122").Value = 23
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 0
$ws.Range("E122").Value = 22
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 1

$ws.Range("A123").Value = "Paraguay"
$ws.Range("B123").Value = 22
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 0
$ws.Range("E123").Value = 21
$ws.Range("F123").Value = 1
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 1

$ws.Range("A124").Value = "Montenegro"
$ws.Range("B124").Value = 22
$ws.Range("C124").Value = 1
$ws.Range("D124").Value = 0
$ws.Range("E124").Value = 21
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 1

$ws.Range("A126").Value = "Mayotte"
$ws.Range("B126").Value = 21
$ws.Range("C126").Value = 10
$ws.Range("D126").Value = 0
$ws.Range("E126").Value = 21
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 0

$ws.Range("A127").Value = "Ruanda"
$ws.Range("B127").Value = 19
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 0
$ws.Range("E127").Value = 19
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 0

$ws.Range("A129").Value = "Guatemala"
$ws.Range("B129").Value = 19
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 0
$ws.Range("E129").Value = 18
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 1

$ws.Range("A130").Value = "Jamaica"
$ws.Range("B130").Value = 19
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 2
$ws.Range("E130").Value = 16
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 1

$ws.Range("A131").Value = "Polinesia Francesa"
$ws.Range("B131").Value = 18
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 0
$ws.Range("E131").Value = 18
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 0

$ws.Range("A132").Value = "Guayana Francesa"
$ws.Range("B132").Value = 18
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 0
$ws.Range("E132").Value = 18
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 0

$ws.Range("A133").Value = "Togo"
$ws.Range("B133").Value = 16
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 0
$ws.Range("E133").Value = 16
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 0

$ws.Range("A134").Value = "Kenia"
$ws.Range("B134").Value = 15
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 0
$ws.Range("E134").Value = 15
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 0

$ws.Range("A135").Value = "Gibraltar"
$ws.Range("B135").Value = 15
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 2
$ws.Range("E135").Value = 13
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 0

$ws.Range("A137").Value = "Kirguistan"
$ws.Range("B137").Value = 14
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 0
$ws.Range("E137").Value = 14
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 0

$ws.Range("A138").Value = "Maldivas"
$ws.Range("B138").Value = 13
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 3
$ws.Range("E138").Value = 10
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0

$ws.Range("A139").Value = "Tanzania"
$ws.Range("B139").Value = 12
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 0
$ws.Range("E139").Value = 12
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 0

$ws.Range("A140").Value = "Madagascar"
$ws.Range("B140").Value = 12
$ws.Range("C140").Value = 9
$ws.Range("D140").Value = 0
$ws.Range("E140").Value = 12
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 0

$ws.Range("A144").Value = "Nueva Caledonia"
$ws.Range("B144").Value = 8
$ws.Range("C144").Value = 4
$ws.Range("D144").Value = 0
$ws.Range("E144").Value = 8
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 0

